$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to be treated as text so Excel does not silently
    # convert numeric-looking strings (e.g. "211.89") into real numbers
    # and does not drop significant trailing zeros.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.403.18"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.573.91"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "211.89"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.05%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - OKB
Set-TextValue $ws.Range("D8") "44.41"
$ws.Range("E8").Value = "  -3.93%  "

# Row 9 - Solana
Set-TextValue $ws.Range("D9") "23.75"
$ws.Range("E9").Value = "  -1.39%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.44%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0588"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12 - TRON
Set-TextValue $ws.Range("D12") "0.0895"
$ws.Range("E12").Value = "  +1.39%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "1.798.46"
$ws.Range("E13").Value = "  +0.17%  "

# Row 14 - WrappedEther
Set-TextValue $ws.Range("D14") "1.565.69"
$ws.Range("E14").Value = "  -0.21%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "28.409.45"
$ws.Range("E16").Value = "  -0.30%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -1.03%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "61.71"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "228.66"
$ws.Range("E19").Value = "  +0.70%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +0.62%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -1.10%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +1.63%  "

# Row 24 - Avalanche
Set-TextValue $ws.Range("D24") "8.96"

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.85%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "150.75"
$ws.Range("E26").Value = "  +0.10%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "14.93"
$ws.Range("E27").Value = "  -0.32%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.31%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  -1.11%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  -0.05%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("D31") "0.0481"
$ws.Range("E31").Value = "  +3.53%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -2.47%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -1.41%  "

# Row 35 - Maker
Set-TextValue $ws.Range("D35") "1.383.47"
$ws.Range("E35").Value = "  -0.68%  "

# Row 36 - TrustWalletToken
Set-TextValue $ws.Range("D36") "1.07"
$ws.Range("E36").Value = "  +4.31%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -2.58%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  -0.37%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  +1.66%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -1.39%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  -2.39%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +3.14%  "

# Row 43 - PaxDollar
$ws.Range("E43").Value = "  +0.02%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  -0.45%  "

# Row 45 - Kaspa
Set-TextValue $ws.Range("D45") "0.0471"
$ws.Range("E45").Value = "  -0.33%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  -4.28%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "62.40"
$ws.Range("E47").Value = "  -0.92%  "

# Row 48 - WEMIXToken
Set-TextValue $ws.Range("D48") "0.919"
$ws.Range("E48").Value = "  -6.20%  "

# Row 49 - RocketPoolETH
Set-TextValue $ws.Range("D49") "1.710.74"
$ws.Range("E49").Value = "  +0.15%  "

# Row 50 - mCoin
$ws.Range("E50").Value = "  +1.81%  "

# Row 51 - Quant
Set-TextValue $ws.Range("D51") "85.54"
$ws.Range("E51").Value = "  -0.51%  "
